# Updates cryptos list values (Price and Volume(1h) columns) per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.720.50"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "3.303.16"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.587"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.27%  "
$ws.Range("D8").Value = "3.302.42"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -4.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -12.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.133"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("D15").Value = "3.829.28"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "3.308.80"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "63.416.36"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "374.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  +4.40%  "
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  -5.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "620.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.71%  "
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.374"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.85%  "
$ws.Range("D40").Value = "0.0₃0724"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("E43").Value = "  -4.87%  "
$ws.Range("D44").Value = "2.870.34"
$ws.Range("E44").Value = "  -2.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0392"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("E48").Value = "  -6.48%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("E51").Value = "  -3.14%  "
